$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 gets a new distinct set of credential values (not the same as row 2)
$ws.Range("A4").Value = "1231@yopmail.com"
$ws.Range("B4").Value = "adm1"

# Move the active selection to B4
$ws.Range("B4").Select()
